$wb = $excel.ActiveWorkbook
$ws6 = $wb.Worksheets.Item("Prabu")
$ws6.Columns.Item(3).ColumnWidth = 30.2109375
